$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 and J1, copying formatting from the existing H1 header cell
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1").PasteSpecial(-4122) | Out-Null
$ws.Range("J1").PasteSpecial(-4122) | Out-Null

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the data values for the new I and J columns
$values = @{
    2  = @(5, 5)
    3  = @(8, 8)
    4  = @(8, 8)
    5  = @(8, 8)
    6  = @(8, 8)
    7  = @(6, 6)
    8  = @(7, 7)
    9  = @(4, 4)
    10 = @(7, 7)
    11 = @(5, 6)
    12 = @(6, 6)
    13 = @(7, 7)
    14 = @(6, 6)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
